$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row label text for the rotated summary rows (A33:A36)
$ws.Range("A33").Value = "TotalNNAvg"
$ws.Range("A34").Value = "Pre2020NNavg"
$ws.Range("A35").Value = "TransNNavg"
$ws.Range("A36").Value = "Post2020NNavg"

# Update numeric metric values for rows 2-36 (columns B-G)
$ws.Range("B2").Value = -0.078
$ws.Range("C2").Value = 1.693
$ws.Range("D2").Value = 7.44
$ws.Range("E2").Value = 2.728
$ws.Range("G2").Value = 0.1
$ws.Range("B3").Value = -0.042
$ws.Range("C3").Value = 1.67
$ws.Range("D3").Value = 0.662
$ws.Range("E3").Value = 0.8139999999999999
$ws.Range("F3").Value = 0.57
$ws.Range("G3").Value = -0.206
$ws.Range("B4").Value = -0.034
$ws.Range("C4").Value = 1.665
$ws.Range("D4").Value = 0.18
$ws.Range("E4").Value = 0.424
$ws.Range("F4").Value = 0.32
$ws.Range("G4").Value = -0.197
$ws.Range("B5").Value = -0.509
$ws.Range("C5").Value = 1.97
$ws.Range("D5").Value = 1.186
$ws.Range("E5").Value = 1.089
$ws.Range("F5").Value = 0.846
$ws.Range("G5").Value = -0.638
$ws.Range("B6").Value = -0.04
$ws.Range("C6").Value = 1.669
$ws.Range("D6").Value = 0.922
$ws.Range("E6").Value = 0.96
$ws.Range("F6").Value = 0.841
$ws.Range("G6").Value = -0.574
$ws.Range("B7").Value = -0.012
$ws.Range("C7").Value = 1.651
$ws.Range("D7").Value = 1.267
$ws.Range("E7").Value = 1.126
$ws.Range("F7").Value = 0.992
$ws.Range("G7").Value = -0.232
$ws.Range("B8").Value = -0.08599999999999999
$ws.Range("C8").Value = 1.698
$ws.Range("D8").Value = 0.372
$ws.Range("E8").Value = 0.61
$ws.Range("F8").Value = 0.483
$ws.Range("G8").Value = 0.041
$ws.Range("B9").Value = -0.191
$ws.Range("C9").Value = 1.766
$ws.Range("D9").Value = 1.389
$ws.Range("E9").Value = 1.179
$ws.Range("F9").Value = 0.964
$ws.Range("G9").Value = 0.142
$ws.Range("B10").Value = -0.005
$ws.Range("C10").Value = 1.646
$ws.Range("D10").Value = 0.901
$ws.Range("E10").Value = 0.949
$ws.Range("F10").Value = 0.839
$ws.Range("G10").Value = 0.014
$ws.Range("B11").Value = 0.001
$ws.Range("C11").Value = 1.642
$ws.Range("D11").Value = 0.63
$ws.Range("E11").Value = 0.794
$ws.Range("F11").Value = 0.679
$ws.Range("G11").Value = 0.261
$ws.Range("B12").Value = -0.004
$ws.Range("C12").Value = 1.645
$ws.Range("D12").Value = 0.303
$ws.Range("E12").Value = 0.55
$ws.Range("F12").Value = 0.43
$ws.Range("G12").Value = 0.045
$ws.Range("B13").Value = -0.07099999999999999
$ws.Range("C13").Value = 1.688
$ws.Range("D13").Value = 0.523
$ws.Range("E13").Value = 0.723
$ws.Range("F13").Value = 0.674
$ws.Range("G13").Value = 0.466
$ws.Range("B14").Value = -0.063
$ws.Range("C14").Value = 1.683
$ws.Range("D14").Value = 0.532
$ws.Range("E14").Value = 0.729
$ws.Range("F14").Value = 0.534
$ws.Range("G14").Value = 0.437
$ws.Range("B15").Value = 0.006
$ws.Range("C15").Value = 1.639
$ws.Range("D15").Value = 0.391
$ws.Range("E15").Value = 0.625
$ws.Range("F15").Value = 0.508
$ws.Range("G15").Value = 0.188
$ws.Range("B16").Value = -0.982
$ws.Range("C16").Value = 1.417
$ws.Range("D16").Value = 0.215
$ws.Range("E16").Value = 0.464
$ws.Range("F16").Value = 0.425
$ws.Range("G16").Value = 0.351
$ws.Range("B17").Value = -2.041
$ws.Range("C17").Value = 1.138
$ws.Range("D17").Value = 0.231
$ws.Range("E17").Value = 0.481
$ws.Range("F17").Value = 0.394
$ws.Range("B18").Value = -7.694
$ws.Range("C18").Value = 1.395
$ws.Range("D18").Value = 0.174
$ws.Range("E18").Value = 0.417
$ws.Range("F18").Value = 0.391
$ws.Range("G18").Value = -1
$ws.Range("B19").Value = -12.55
$ws.Range("C19").Value = 1.616
$ws.Range("D19").Value = 0.6840000000000001
$ws.Range("E19").Value = 0.827
$ws.Range("F19").Value = 0.788
$ws.Range("G19").Value = -1
$ws.Range("B20").Value = -9.769
$ws.Range("C20").Value = 1.49
$ws.Range("D20").Value = 0.336
$ws.Range("E20").Value = 0.58
$ws.Range("F20").Value = 0.556
$ws.Range("B21").Value = -64.40900000000001
$ws.Range("C21").Value = 3.973
$ws.Range("D21").Value = 1.768
$ws.Range("E21").Value = 1.33
$ws.Range("F21").Value = 1.321
$ws.Range("G21").Value = 1
$ws.Range("B22").Value = -10.36
$ws.Range("C22").Value = 1.516
$ws.Range("D22").Value = 0.386
$ws.Range("E22").Value = 0.621
$ws.Range("F22").Value = 0.591
$ws.Range("G22").Value = -1
$ws.Range("B23").Value = -8.734999999999999
$ws.Range("C23").Value = 1.442
$ws.Range("D23").Value = 1.537
$ws.Range("E23").Value = 1.24
$ws.Range("F23").Value = 1.147
$ws.Range("B24").Value = -0.2
$ws.Range("C24").Value = 1.055
$ws.Range("D24").Value = 0.786
$ws.Range("E24").Value = 0.887
$ws.Range("F24").Value = 0.8100000000000001
$ws.Range("B25").Value = -6.954
$ws.Range("C25").Value = 1.362
$ws.Range("D25").Value = 0.5649999999999999
$ws.Range("E25").Value = 0.752
$ws.Range("F25").Value = 0.718
$ws.Range("B26").Value = -0.316
$ws.Range("C26").Value = 1.06
$ws.Range("D26").Value = 0.014
$ws.Range("E26").Value = 0.118
$ws.Range("F26").Value = 0.113
$ws.Range("G26").Value = 1
$ws.Range("B27").Value = -0.574
$ws.Range("C27").Value = 1.072
$ws.Range("D27").Value = 1.991
$ws.Range("E27").Value = 1.411
$ws.Range("G27").Value = 1
$ws.Range("B28").Value = -2.348
$ws.Range("C28").Value = 1.152
$ws.Range("D28").Value = 0.516
$ws.Range("E28").Value = 0.718
$ws.Range("F28").Value = 0.625
$ws.Range("G28").Value = 1
$ws.Range("B29").Value = -6979.821
$ws.Range("C29").Value = 318.31
$ws.Range("D29").Value = 8.167
$ws.Range("E29").Value = 2.858
$ws.Range("F29").Value = 2.857
$ws.Range("G29").Value = -1
$ws.Range("B30").Value = -0.049
$ws.Range("C30").Value = 1.048
$ws.Range("D30").Value = 0.08
$ws.Range("E30").Value = 0.283
$ws.Range("G30").Value = 1
$ws.Range("B31").Value = -27.687
$ws.Range("C31").Value = 2.304
$ws.Range("D31").Value = 2.551
$ws.Range("E31").Value = 1.597
$ws.Range("F31").Value = 1.564
$ws.Range("G31").Value = -1
$ws.Range("B32").Value = -10.554
$ws.Range("C32").Value = 1.525
$ws.Range("D32").Value = 1.41
$ws.Range("E32").Value = 1.187
$ws.Range("F32").Value = 1.123
$ws.Range("B33").Value = -52.86607878151261
$ws.Range("C33").Value = 3.924334033613445
$ws.Range("D33").Value = 0.9206848739495799
$ws.Range("E33").Value = 0.8078455882352942
$ws.Range("F33").Value = 0.6905388655462186
$ws.Range("G33").Value = 0.08210084033613443
$ws.Range("B34").Value = -0.08057142857142856
$ws.Range("C34").Value = 1.694642857142857
$ws.Range("D34").Value = 1.192714285714286
$ws.Range("E34").Value = 0.9500000000000001
$ws.Range("F34").Value = 0.7814285714285715
$ws.Range("G34").Value = -0.01092857142857145
$ws.Range("B35").Value = -0.982
$ws.Range("C35").Value = 1.417
$ws.Range("D35").Value = 0.215
$ws.Range("E35").Value = 0.464
$ws.Range("F35").Value = 0.425
$ws.Range("G35").Value = 0.351
$ws.Range("B36").Value = -446.5038125
$ws.Range("C36").Value = 21.341125
$ws.Range("D36").Value = 1.32475
$ws.Range("E36").Value = 0.9566875
$ws.Range("F36").Value = 0.8999374999999999
$ws.Range("G36").Value = -0.125

Write-Host "Updated full revin NN tuning metrics to new laptop values"
